# Adds a new "2022-Q1" worksheet (with fund-holding detail rows) between the
# existing "2021-Q4" sheet and the "总计" (totals) summary sheet, and inserts
# a corresponding "2022-Q1" row at the top of the "总计" sheet's data table.

$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item(1)      # "2021-Q4" - used as a formatting template

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q1" sheet right after "2021-Q4"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $src)
$ws2.Name = "2022-Q1"

# Re-resolve the totals sheet by name AFTER inserting - grabbing it earlier
# (by index/reference) ends up stale once the sheet collection is mutated.
$totals = $wb.Worksheets.Item("总计")

# Header row: copy style from the template sheet, then overwrite the text
$src.Range("B1:H1").Copy()
$ws2.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats

$ws2.Range("B1").Value = "基金代码"
$ws2.Range("C1").Value = "基金名称"
$ws2.Range("D1").Value = "基金规模"
$ws2.Range("E1").Value = "股票总仓位"
$ws2.Range("F1").Value = "仓位占比"
$ws2.Range("G1").Value = "持有市值(亿元)"
$ws2.Range("H1").Value = "仓位排名"

# Data rows: col A = index (number, bold/bordered style), B/C/D/E/F/G = text,
# H = plain number
$data2 = @(
    ,@(0, '501079', '大成科创主题 3 年封闭运作灵活配置混合', '17.69', '79.13', '4.62', '0.8173', 4)
    ,@(1, '011637', '广发沪港深价值成长混合型证券投资基金A', '12.44', '92.96', '5.37', '0.6680', 7)
    ,@(2, '010452', '广发瑞福精选混合A', '16.29', '78.69', '2.79', '0.4545', 8)
    ,@(3, '012473', '大成成长回报六个月持有期混合型证券投资基金A', '8.97', '71.30', '4.67', '0.4189', 2)
    ,@(4, '010371', '大成成长进取混合A', '5.55', '80.17', '4.69', '0.2603', 4)
    ,@(5, '005743', '长安裕隆灵活配置混合A', '5.68', '91.64', '3.43', '0.1948', 9)
    ,@(6, '013513', '长安先进制造混合A', '2.95', '91.34', '4.61', '0.1360', 6)
    ,@(7, '005341', '长安裕泰灵活配置混合A', '1.73', '91.20', '4.73', '0.0818', 5)
    ,@(8, '010372', '大成成长进取混合C', '1.71', '80.17', '4.69', '0.0802', 4)
    ,@(9, '005744', '长安裕隆灵活配置混合C', '2.13', '91.64', '3.43', '0.0731', 9)
    ,@(10, '010453', '广发瑞福精选混合C', '2.42', '78.69', '2.79', '0.0675', 8)
    ,@(11, '005049', '长安鑫旺价值灵活配置混合A', '1.70', '90.83', '3.60', '0.0612', 10)
    ,@(12, '005342', '长安裕泰灵活配置混合C', '1.22', '91.20', '4.73', '0.0577', 5)
    ,@(13, '005050', '长安鑫旺价值灵活配置混合C', '1.08', '90.83', '3.60', '0.0389', 10)
    ,@(14, '013514', '长安先进制造混合C', '0.52', '91.34', '4.61', '0.0240', 6)
    ,@(15, '012474', '大成成长回报六个月持有期混合型证券投资基金C', '0.43', '71.30', '4.67', '0.0201', 2)
    ,@(16, '011638', '广发沪港深价值成长混合型证券投资基金C', '0.33', '92.96', '5.37', '0.0177', 7)
    ,@(17, '710301', '富安达增强收益债券A', '0.61', '20.20', '2.35', '0.0143', 3)
    ,@(18, '740001', '长安宏观策略混合', '0.16', '71.93', '6.92', '0.0111', 1)
    ,@(19, '710302', '富安达增强收益债券C', '0.26', '20.20', '2.35', '0.0061', 3)
)

$r = 2
foreach ($row in $data2) {
    # Force B..G to be stored as *text* (they hold numeric-looking strings in
    # the source data, e.g. "17.69"), not auto-converted to numbers.
    $ws2.Cells.Item($r, 2).NumberFormat = "@"
    $ws2.Cells.Item($r, 3).NumberFormat = "@"
    $ws2.Cells.Item($r, 4).NumberFormat = "@"
    $ws2.Cells.Item($r, 5).NumberFormat = "@"
    $ws2.Cells.Item($r, 6).NumberFormat = "@"
    $ws2.Cells.Item($r, 7).NumberFormat = "@"

    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $ws2.Cells.Item($r, 6).Value = $row[5]
    $ws2.Cells.Item($r, 7).Value = $row[6]
    $ws2.Cells.Item($r, 8).Value = $row[7]

    # Drop the temporary text number-format again so the cells end up with no
    # explicit style, matching the plain (unstyled) data cells elsewhere.
    $rowTextRange = $ws2.Range($ws2.Cells.Item($r, 2), $ws2.Cells.Item($r, 7))
    $rowTextRange.ClearFormats()

    $r = $r + 1
}
$lastDataRow = $r - 1

# Column A (the row-index column) keeps the bold/bordered/centred style used
# throughout the workbook - copy it from the template sheet.
$src.Range("A2").Copy()
$colARange = $ws2.Range($ws2.Cells.Item(2, 1), $ws2.Cells.Item($lastDataRow, 1))
$colARange.PasteSpecial(-4122)  # xlPasteFormats

# ---------------------------------------------------------------------------
# 2) Add a "2022-Q1" row at the top of the "总计" sheet's data table
# ---------------------------------------------------------------------------
$totals.Rows.Item(2).Insert()

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q1"
$totals.Range("C2").Value = 20
$totals.Range("D2").Value = 3.5
$totals.Range("B2:D2").ClearFormats()

# Re-apply the index-column style to A2, and renumber the (shifted) old row
$totals.Range("A3").Copy()
$totals.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$totals.Range("A3").Value = 1
